$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F17").Value = 14571
$ws1.Range("F20").Value = 577
$ws1.Range("F22").Value = 148
$ws1.Range("F23").Value = 542
$ws1.Range("F25").Value = 109
$ws1.Range("F26").Value = 48
$ws1.Range("F31").Value = 38

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5747

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F28").Value = 577
$ws4.Range("F30").Value = 148
$ws4.Range("F31").Value = 542
$ws4.Range("F33").Value = 109
$ws4.Range("F35").Value = 48
$ws4.Range("F49").Value = 38
